{"js": "// Replace the two Word field constructs (\" m:userdoc 'zone1' \" and\n// \" m:enduserdoc \") with plain literal text runs using the same\n// brace-delimited token syntax (\"{m:userdoc 'zone1'}\" / \"{m:enduserdoc}\"),\n// mirroring the TokenIteratorFieldRewriterSplit parser change.\n//\n// Original runs (zone1 paragraph):\n//   fldChar(begin), instrText(\" \"), instrText(\"m\"),\n//   instrText(\":userdoc 'zone1'\"), instrText(\" \"), fldChar(end)\n// Target runs:\n//   t(\"{\"), t(\"m\"), t(\":userdoc 'zone1'\"), t(\"}\")\n//\n// Original runs (enduserdoc paragraph):\n//   fldChar(begin), instrText(\" m:\"), bookmarkStart/End,\n//   instrText(\"enduserdoc \"), fldChar(end)\n// Target runs:\n//   t(\"{m:\"), bookmarkStart/End (unchanged), t(\"enduserdoc}\")\n\nconst zone1Xml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{</w:t></w:r>\n            <w:r><w:t>m</w:t></w:r>\n            <w:r><w:t>:userdoc 'zone1'</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst endDocXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{m:</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n// Locate the two field paragraphs by inspecting their OOXML content (robust\n// against paragraph-index drift) rather than assuming fixed indices.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet zone1ParagraphIndex = -1;\nlet endDocParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const ooxml = paragraphs.items[i].getOoxml();\n  await context.sync();\n  const xml = ooxml.value;\n  if (xml.indexOf(\"enduserdoc\") !== -1) {\n    endDocParagraphIndex = i;\n  } else if (xml.indexOf(\"userdoc\") !== -1) {\n    zone1ParagraphIndex = i;\n  }\n}\n\nif (zone1ParagraphIndex !== -1) {\n  const zone1Range = paragraphs.items[zone1ParagraphIndex].getRange(\"Content\");\n  zone1Range.insertOoxml(zone1Xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nif (endDocParagraphIndex !== -1) {\n  const endRange = paragraphs.items[endDocParagraphIndex].getRange(\"Content\");\n  endRange.insertOoxml(endDocXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the two Word field constructs (\" m:userdoc 'zone1' \" and\n# \" m:enduserdoc \") with plain literal text runs using the same\n# brace-delimited token syntax (\"{m:userdoc 'zone1'}\" / \"{m:enduserdoc}\"),\n# mirroring the TokenIteratorFieldRewriterSplit parser change.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexForPosition($doc, $pos) {\n  for ($j = 1; $j -le $doc.Paragraphs.Count; $j++) {\n    $pr = $doc.Paragraphs.Item($j).Range\n    if ($pos -ge $pr.Start -and $pos -lt $pr.End) {\n      return $j\n    }\n  }\n  return -1\n}\n\n# Target runs: t(\"{\"), t(\"m\"), t(\":userdoc 'zone1'\"), t(\"}\")\n$zone1Xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ''zone1''</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Target runs: t(\"{m:\"), bookmarkStart/End (kept as-is), t(\"enduserdoc}\")\n$enddocXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Locate the paragraphs that host each field by inspecting the field codes\n# directly (robust against paragraph-index drift); $d.Fields exposes correct\n# Code.Start/End offsets that we map back to a paragraph index.\n$zone1ParaIdx = -1\n$enddocParaIdx = -1\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n  $f = $d.Fields.Item($i)\n  $code = $f.Code.Text\n  $idx = Find-ParagraphIndexForPosition $d $f.Code.Start\n  if ($code -match \"userdoc\" -and $code -notmatch \"enduserdoc\") {\n    $zone1ParaIdx = $idx\n  } elseif ($code -match \"enduserdoc\") {\n    $enddocParaIdx = $idx\n  }\n}\n\nif ($zone1ParaIdx -gt 0) {\n  $d.Paragraphs.Item($zone1ParaIdx).Range.InsertXML($zone1Xml)\n}\n\n# Re-fetch the document/paragraphs after the first mutation before touching\n# the second field's paragraph.\n$d2 = $word.ActiveDocument\nif ($enddocParaIdx -gt 0) {\n  $d2.Paragraphs.Item($enddocParaIdx).Range.InsertXML($enddocXml)\n}\n"}
